$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) keeps text formatting so values like "1.00" or "0.0841"
# are not silently converted into numbers by Excel when assigned.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "40.889.53"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").Value = "2.387.29"
$ws.Range("E3").Value = "  -3.71%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "313.79"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("D6").Value = "88.20"
$ws.Range("E6").Value = "  -5.72%  "
$ws.Range("D7").Value = "0.531"
$ws.Range("E7").Value = "  -3.93%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  -3.96%  "
$ws.Range("D10").Value = "0.0841"
$ws.Range("D11").Value = "31.10"
$ws.Range("E11").Value = "  -6.85%  "
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").Value = "2.756.83"
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("D14").Value = "6.61"
$ws.Range("E14").Value = "  -5.14%  "
$ws.Range("D15").Value = "15.24"
$ws.Range("E15").Value = "  -3.08%  "
$ws.Range("D16").Value = "2.392.84"
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("D17").Value = "0.768"
$ws.Range("E17").Value = "  -4.01%  "
$ws.Range("D18").Value = "40.846.13"
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("E19").Value = "  -3.86%  "
$ws.Range("E20").Value = "  -4.61%  "
$ws.Range("D21").Value = "69.07"
$ws.Range("E21").Value = "  -3.09%  "
$ws.Range("D22").Value = "11.00"
$ws.Range("E22").Value = "  -3.59%  "
$ws.Range("D23").Value = "234.77"
$ws.Range("E23").Value = "  -3.11%  "
$ws.Range("E24").Value = "  -3.93%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -6.93%  "
$ws.Range("D27").Value = "23.98"
$ws.Range("E27").Value = "  -5.38%  "
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").Value = "9.42"
$ws.Range("E29").Value = "  -3.98%  "
$ws.Range("D30").Value = "34.06"
$ws.Range("E30").Value = "  -8.58%  "
$ws.Range("D31").Value = "153.42"
$ws.Range("E31").Value = "  -2.72%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "5.29"
$ws.Range("E32").Value = "  -4.31%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "0.0736"
$ws.Range("E34").Value = "  -4.29%  "
$ws.Range("E35").Value = "  -5.09%  "
$ws.Range("E36").Value = "  -2.24%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "16.24"
$ws.Range("E37").Value = "  -7.59%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "2.82"
$ws.Range("E38").Value = "  -4.21%  "
$ws.Range("E39").Value = "  -4.32%  "
$ws.Range("D40").Value = "1.73"
$ws.Range("E40").Value = "  -8.03%  "
$ws.Range("D41").Value = "3.86"
$ws.Range("E41").Value = "  -4.20%  "
$ws.Range("E42").Value = "  -7.27%  "
$ws.Range("D43").Value = "1.971.23"
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("D44").Value = "0.0272"
$ws.Range("E44").Value = "  -4.65%  "
$ws.Range("D45").Value = "17.96"
$ws.Range("E45").Value = "  -6.91%  "
$ws.Range("D46").Value = "9.66"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("D47").Value = "2.75"
$ws.Range("E47").Value = "  -8.16%  "
$ws.Range("D48").Value = "2.620.55"
$ws.Range("E48").Value = "  -3.60%  "
$ws.Range("D49").Value = "94.11"
$ws.Range("E49").Value = "  -4.52%  "
$ws.Range("D50").Value = "72.91"
$ws.Range("E50").Value = "  -4.61%  "
$ws.Range("D51").Value = "51.14"
$ws.Range("E51").Value = "  -2.80%  "
